# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" for the b2121808-9ac1-4fd6-a3d6-22fe1d966b9c
# row (row 7) on the zh-cn and de-de language sheets, and refreshes the
# corresponding "Latest HO Xliff Generate Date" on the Overview sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-31 14:52:51"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-31 14:52:47"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-31 14:52:51"
